# cryptos.xlsx data refresh (Thu Mar 30 19:51:02 UTC 2023, GitHub Actions)
#
# The Price (D) / Volume(1h) (E) columns hold plain-text values (e.g.
# "316.24", "  +0.06%  ") rather than numbers, so every assignment below
# is apostrophe-prefixed to force text entry (stops Excel from silently
# reinterpreting number-looking strings as floating point values), and the
# cell style is put back to "Normal" right after so no stray quote-prefix
# formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.010.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.24%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.781.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.36%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'316.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.06%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.07%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.5379"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.36%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3769"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.25%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07442"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.98%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'41.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.02%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -2.57%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +0.03%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'20.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.21%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.082"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.76%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'7.223"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.58%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'1.775.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.17%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -4.44%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -1.60%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.06434"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.14%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.06%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'17.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.90%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.871"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.20%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'28.026.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.13%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.06%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.086"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.96%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'155.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.78%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'20.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.37%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.978.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.43%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.279"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.02%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'119.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.06%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -1.77%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.1053"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.32%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.69%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'5.513"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.05%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.2248"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.24%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.06448"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.06%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.02282"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.70%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'4.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.50%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'8.409"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.76%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D42").Value = "'1.438"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.04%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.176"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.44%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.9995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.02%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'13.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.46%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'3.669"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.48%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5740"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.06%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'126.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.72%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.185"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.40%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.921"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.23%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.06791"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.66%  "
$ws.Range("E51").Style = "Normal"

# Rows 40/41: Aptos and TheSandbox swapped positions in the ranking
$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.6134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.46%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'Aptos"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'11.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.93%  "
$ws.Range("E41").Style = "Normal"
